# Controle Transferencia.xlsx - automatic worksheet update
# - Row 8 (A8) switches from the "date only" style to the "date + time" style
#   (it is no longer the last/"in progress" row).
# - A new row 9 is appended, reusing the "date only" style that A8 used to
#   have, for a new truck "asd5678" / conferente "solo" that has just
#   started (only Data/Placa/Conferente are filled in; every other column
#   is an empty, but present, text cell - mirroring the other "in progress"
#   rows such as row 3 and row 5).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember A8's current ("date only") number format before we change it -
# the new row 9 will reuse it.
$dateOnlyFormat = $ws.Range("A8").NumberFormat

# A8 now gets the "date + time" format/style (same style already used by
# A2..A7).
$ws.Range("A8").NumberFormat = $ws.Range("A2").NumberFormat

# Append row 9.
$ws.Range("A9").Value = 45856
$ws.Range("A9").NumberFormat = $dateOnlyFormat

$ws.Range("B9").Value = "asd5678"
$ws.Range("C9").Value = "solo"

# D9:W9 stay empty (no data yet for this new truck) but, like the rest of
# the sheet, are real empty text cells rather than untouched/missing ones.
# Writing a plain "" clears/removes a cell in this engine (matches real
# Excel semantics), so use the classic force-text leading apostrophe to get
# an empty *text* cell, then drop back to the Normal style so no stray
# quote-prefix formatting is left on the cell.
for ($c = 4; $c -le 23; $c++) {
    $cell = $ws.Cells.Item(9, $c)
    $cell.Value = "'"
    $cell.Style = "Normal"
}
